$wb = $excel.ActiveWorkbook

# --- Sheet "ip_address_list": add row 3 ---
$ws1 = $wb.Worksheets.Item("ip_address_list")
$row1 = $ws1.Range("A3:C3")
$row1.NumberFormat = "@"
$ws1.Range("A3").Value = "515"
$ws1.Range("B3").Value = "10.9.250.241"
$ws1.Range("C3").Value = "255.255.255.0"
$row1.ClearFormats()
$ws1.Range("E3").Value = 0

# --- Sheet "disk_list": add row 5 ---
$ws3 = $wb.Worksheets.Item("disk_list")
$row3 = $ws3.Range("A5:E5")
$row3.NumberFormat = "@"
$ws3.Range("A5").Value = "515"
$ws3.Range("B5").Value = "Z"
$ws3.Range("C5").Value = "\\10.9.250.100\"
$ws3.Range("D5").Value = "spravce"
$ws3.Range("E5").Value = "Jhv*2708"
$row3.ClearFormats()
